$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 2663.3333
$ws.Range("J17").Value = 2663.3333
$ws.Range("L17").Value = 7989.999899999999
$ws.Range("N17").Value = -8325.999899999999
# Row 33
$ws.Range("H33").Value = 25000776
$ws.Range("I33").Value = 31250688
$ws.Range("J33").Value = 1125
$ws.Range("K33").Value = 31250688
$ws.Range("L33").Value = 1125
$ws.Range("M33").Value = -31250459
$ws.Range("N33").Value = -1583
# Row 62
$ws.Range("H62").Value = 25999.4
$ws.Range("I62").Value = 25999.4
$ws.Range("K62").Value = 25999.4
$ws.Range("M62").Value = -25375.4
# Row 65
$ws.Range("H65").Value = 25999.4
$ws.Range("I65").Value = 25999.4
$ws.Range("K65").Value = 129997
$ws.Range("M65").Value = -126877
# Row 74
$ws.Range("H74").Value = 9205.25
$ws.Range("I74").Value = 9085.6
$ws.Range("K74").Value = 9085.6
$ws.Range("M74").Value = -8149.6
# Row 77
$ws.Range("H77").Value = 9205.25
$ws.Range("I77").Value = 9085.6
$ws.Range("K77").Value = 45428
$ws.Range("M77").Value = -40748
# Row 98
$ws.Range("H98").Value = 622.53845
$ws.Range("I98").Value = 622.53845
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 622.53845
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 875.46155
$ws.Range("N98").ClearContents()
# Row 107
$ws.Range("H107").Value = 508.5
$ws.Range("I107").Value = 539.61536
$ws.Range("J107").Value = 373.66666
$ws.Range("K107").Value = 539.61536
$ws.Range("L107").Value = 373.66666
$ws.Range("M107").Value = 1380.38464
$ws.Range("N107").Value = -4213.66666
# Row 116
$ws.Range("H116").Value = 3766.7144
$ws.Range("I116").Value = 4274.4
$ws.Range("J116").Value = 2497.5
$ws.Range("K116").Value = 4274.4
$ws.Range("L116").Value = 2497.5
$ws.Range("M116").Value = -832.3999999999996
$ws.Range("N116").Value = -9381.5
# Row 122
$ws.Range("H122").Value = 622.53845
$ws.Range("I122").Value = 622.53845
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 1867.61535
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 582.38465
$ws.Range("N122").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2309.75
$ws.Range("I2").Value = 2102.8572
$ws.Range("J2").Value = 2599.4
$ws.Range("K2").Value = 2102.8572
$ws.Range("L2").Value = 2599.4
$ws.Range("M2").Value = -1989.8572
$ws.Range("N2").Value = -2825.4
# Row 32
$ws.Range("H32").Value = 5635.174
$ws.Range("I32").Value = 4271.3022
$ws.Range("J32").Value = 25184
$ws.Range("K32").Value = 4271.3022
$ws.Range("L32").Value = 25184
$ws.Range("M32").Value = -3984.3022
$ws.Range("N32").Value = -25758
# Row 45
$ws.Range("H45").Value = 706
$ws.Range("I45").Value = 706
$ws.Range("K45").Value = 706
$ws.Range("M45").Value = -329
# Row 116
$ws.Range("H116").Value = 2309.75
$ws.Range("I116").Value = 2102.8572
$ws.Range("J116").Value = 2599.4
$ws.Range("K116").Value = 2102.8572
$ws.Range("L116").Value = 2599.4
$ws.Range("M116").Value = 191.1428000000001
$ws.Range("N116").Value = -7187.4
# Row 122
$ws.Range("H122").Value = 6623.9165
$ws.Range("I122").Value = 6635.273
$ws.Range("J122").Value = 6499
$ws.Range("K122").Value = 19905.819
$ws.Range("L122").Value = 19497
$ws.Range("M122").Value = -17455.819
$ws.Range("N122").Value = -24397

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2309.75
$ws.Range("I3").Value = 2102.8572
$ws.Range("J3").Value = 2599.4
$ws.Range("K3").Value = 2102.8572
$ws.Range("L3").Value = 2599.4
$ws.Range("M3").Value = -1988.8572
$ws.Range("N3").Value = -2827.4
# Row 107
$ws.Range("H107").Value = 2825.842
$ws.Range("I107").Value = 2871.7222
$ws.Range("K107").Value = 2871.7222
$ws.Range("M107").Value = -951.7222000000002

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 399.25
$ws.Range("I16").Value = 366
$ws.Range("K16").Value = 366
$ws.Range("M16").Value = -79
# Row 55
$ws.Range("H55").Value = 28995.334
$ws.Range("I55").Value = 36999
$ws.Range("J55").Value = 24993.5
$ws.Range("K55").Value = 36999
$ws.Range("L55").Value = 24993.5
$ws.Range("M55").Value = -36684
$ws.Range("N55").Value = -25623.5
# Row 105
$ws.Range("H105").Value = 1206.4
$ws.Range("I105").Value = 1206.4
$ws.Range("K105").Value = 1206.4
$ws.Range("M105").Value = 540.5999999999999
# Row 107
$ws.Range("H107").Value = 291.16666
$ws.Range("I107").Value = 289.33334
$ws.Range("K107").Value = 289.33334
$ws.Range("M107").Value = 1630.66666
# Row 113
$ws.Range("H113").Value = 399.25
$ws.Range("I113").Value = 366
$ws.Range("K113").Value = 366
$ws.Range("M113").Value = 1804
# Row 141
$ws.Range("H141").Value = 372809
$ws.Range("J141").Value = 372809
$ws.Range("L141").Value = 372809
$ws.Range("N141").Value = -383169

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 8
$ws.Range("H8").Value = 196
$ws.Range("I8").Value = 196
$ws.Range("K8").Value = 588
$ws.Range("M8").Value = -449
# Row 68
$ws.Range("H68").Value = 899.75
$ws.Range("I68").Value = 943.2
$ws.Range("K68").Value = 2829.6
$ws.Range("M68").Value = -2018.6
# Row 71
$ws.Range("H71").Value = 899.75
$ws.Range("I71").Value = 943.2
$ws.Range("K71").Value = 8488.800000000001
$ws.Range("M71").Value = -4432.800000000001
# Row 107
$ws.Range("H107").Value = 369.52
$ws.Range("J107").Value = 361.8
$ws.Range("L107").Value = 1085.4
$ws.Range("N107").Value = -4925.4
# Row 112
$ws.Range("H112").Value = 8615.579
$ws.Range("J112").Value = 9959.866
$ws.Range("L112").Value = 29879.598
$ws.Range("N112").Value = -32095.598
# Row 117
$ws.Range("H117").Value = 840.63635
$ws.Range("I117").Value = 779.44446
$ws.Range("K117").Value = 2338.33338
$ws.Range("M117").Value = 1103.66662

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 31
$ws.Range("H31").Value = 669.6
$ws.Range("I31").Value = 669.6
$ws.Range("K31").Value = 669.6
$ws.Range("M31").Value = -377.6
# Row 37
$ws.Range("H37").Value = 669.6
$ws.Range("I37").Value = 669.6
$ws.Range("K37").Value = 669.6
$ws.Range("M37").Value = -392.6
# Row 102
$ws.Range("H102").Value = 1733.2307
$ws.Range("I102").Value = 1733.2307
$ws.Range("K102").Value = 1733.2307
$ws.Range("M102").Value = -111.2307000000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 166667650
$ws.Range("I22").Value = 954
$ws.Range("K22").Value = 954
$ws.Range("M22").Value = -659
# Row 27
$ws.Range("H27").Value = 166667650
$ws.Range("I27").Value = 954
$ws.Range("K27").Value = 954
$ws.Range("M27").Value = -847
# Row 40
$ws.Range("H40").Value = 37043884
$ws.Range("I40").Value = 58829468
$ws.Range("K40").Value = 58829468
$ws.Range("M40").Value = -58829332
# Row 46
$ws.Range("H46").Value = 3230
$ws.Range("I46").Value = 3230
$ws.Range("K46").Value = 3230
$ws.Range("M46").Value = -3042
# Row 61
$ws.Range("H61").Value = 3999.5334
$ws.Range("I61").Value = 4111
$ws.Range("K61").Value = 4111
$ws.Range("M61").Value = -3909
# Row 113
$ws.Range("H113").Value = 3999.5334
$ws.Range("I113").Value = 4111
$ws.Range("K113").Value = 4111
$ws.Range("M113").Value = -1941
# Row 132
$ws.Range("H132").Value = 7468.4
$ws.Range("J132").Value = 26099.5
$ws.Range("L132").Value = 78298.5
$ws.Range("N132").Value = -83358.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 675.0769
$ws.Range("I113").Value = 660
$ws.Range("K113").Value = 1980
$ws.Range("M113").Value = 190
# Row 132
$ws.Range("H132").Value = 100001300
$ws.Range("I132").Value = 1212.7142
$ws.Range("J132").Value = 333334850
$ws.Range("K132").Value = 3638.1426
$ws.Range("L132").Value = 1000004550
$ws.Range("M132").Value = -1108.1426
$ws.Range("N132").Value = -1000009610

